$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 11.122774
$ws.Cells.Item(2, 8).Value = 33.368322
$ws.Cells.Item(2, 9).Value = 0.2449652610853511
$ws.Cells.Item(2, 10).Value = 0.2449652610853511
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 161.7750676666667
$ws.Cells.Item(2, 14).Value = 485.325203
$ws.Cells.Item(2, 15).Value = 0.9790864123038654
$ws.Cells.Item(2, 16).Value = 0.9790864123038654
$ws.Cells.Item(2, 17).Value = 1799.387516491041
$ws.Cells.Item(2, 18).Value = 16194.48764841936
$ws.Cells.Item(2, 19).Value = 0.2398421586151361
$ws.Cells.Item(2, 20).Value = 0.2398421586151361

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 11.122774
$ws.Cells.Item(3, 8).Value = 33.368322
$ws.Cells.Item(3, 9).Value = 0.2449652610853511
$ws.Cells.Item(3, 10).Value = 0.2449652610853511
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.67894
$ws.Cells.Item(3, 14).Value = 2.03682
$ws.Cells.Item(3, 15).Value = 0.004109044356199978
$ws.Cells.Item(3, 16).Value = 0.004109044356199979
$ws.Cells.Item(3, 17).Value = 7.551696179559999
$ws.Cells.Item(3, 18).Value = 67.96526561604
$ws.Cells.Item(3, 19).Value = 0.001006573123527816
$ws.Cells.Item(3, 20).Value = 0.001006573123527816

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 11.122774
$ws.Cells.Item(4, 8).Value = 33.368322
$ws.Cells.Item(4, 9).Value = 0.2449652610853511
$ws.Cells.Item(4, 10).Value = 0.2449652610853511
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 1.763201333333333
$ws.Cells.Item(4, 14).Value = 5.289604
$ws.Cells.Item(4, 15).Value = 0.01067115280816804
$ws.Cells.Item(4, 16).Value = 0.01067115280816804
$ws.Cells.Item(4, 17).Value = 19.61168994716533
$ws.Cells.Item(4, 18).Value = 176.505209524488
$ws.Cells.Item(4, 19).Value = 0.002614061733734562
$ws.Cells.Item(4, 20).Value = 0.002614061733734562

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 11.122774
$ws.Cells.Item(5, 8).Value = 33.368322
$ws.Cells.Item(5, 9).Value = 0.2449652610853511
$ws.Cells.Item(5, 10).Value = 0.2449652610853511
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 1.013424
$ws.Cells.Item(5, 14).Value = 3.040272
$ws.Cells.Item(5, 15).Value = 0.006133390531766587
$ws.Cells.Item(5, 16).Value = 0.006133390531766588
$ws.Cells.Item(5, 17).Value = 11.272086118176
$ws.Cells.Item(5, 18).Value = 101.448775063584
$ws.Cells.Item(5, 19).Value = 0.001502467612952622
$ws.Cells.Item(5, 20).Value = 0.001502467612952623

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 18.220714
$ws.Cells.Item(6, 8).Value = 54.662142
$ws.Cells.Item(6, 9).Value = 0.4012885600454987
$ws.Cells.Item(6, 10).Value = 0.4012885600454988
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 161.7750676666667
$ws.Cells.Item(6, 14).Value = 485.325203
$ws.Cells.Item(6, 15).Value = 0.9790864123038654
$ws.Cells.Item(6, 16).Value = 0.9790864123038654
$ws.Cells.Item(6, 17).Value = 2947.657240284981
$ws.Cells.Item(6, 18).Value = 26528.91516256483
$ws.Cells.Item(6, 19).Value = 0.3928961765535316
$ws.Cells.Item(6, 20).Value = 0.3928961765535317

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 18.220714
$ws.Cells.Item(7, 8).Value = 54.662142
$ws.Cells.Item(7, 9).Value = 0.4012885600454987
$ws.Cells.Item(7, 10).Value = 0.4012885600454988
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 0.67894
$ws.Cells.Item(7, 14).Value = 2.03682
$ws.Cells.Item(7, 15).Value = 0.004109044356199978
$ws.Cells.Item(7, 16).Value = 0.004109044356199979
$ws.Cells.Item(7, 17).Value = 12.37077156316
$ws.Cells.Item(7, 18).Value = 111.33694406844
$ws.Cells.Item(7, 19).Value = 0.001648912492862573
$ws.Cells.Item(7, 20).Value = 0.001648912492862573

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 18.220714
$ws.Cells.Item(8, 8).Value = 54.662142
$ws.Cells.Item(8, 9).Value = 0.4012885600454987
$ws.Cells.Item(8, 10).Value = 0.4012885600454988
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 1.763201333333333
$ws.Cells.Item(8, 14).Value = 5.289604
$ws.Cells.Item(8, 15).Value = 0.01067115280816804
$ws.Cells.Item(8, 16).Value = 0.01067115280816804
$ws.Cells.Item(8, 17).Value = 32.12678721908533
$ws.Cells.Item(8, 18).Value = 289.141084971768
$ws.Cells.Item(8, 19).Value = 0.004282211544415234
$ws.Cells.Item(8, 20).Value = 0.004282211544415234

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 18.220714
$ws.Cells.Item(9, 8).Value = 54.662142
$ws.Cells.Item(9, 9).Value = 0.4012885600454987
$ws.Cells.Item(9, 10).Value = 0.4012885600454988
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 1.013424
$ws.Cells.Item(9, 14).Value = 3.040272
$ws.Cells.Item(9, 15).Value = 0.006133390531766587
$ws.Cells.Item(9, 16).Value = 0.006133390531766588
$ws.Cells.Item(9, 17).Value = 18.465308864736
$ws.Cells.Item(9, 18).Value = 166.187779782624
$ws.Cells.Item(9, 19).Value = 0.00246125945468931
$ws.Cells.Item(9, 20).Value = 0.00246125945468931

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 0.1189986666666667
$ws.Cells.Item(10, 8).Value = 0.356996
$ws.Cells.Item(10, 9).Value = 0.002620797603979787
$ws.Cells.Item(10, 10).Value = 0.002620797603979787
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 161.7750676666667
$ws.Cells.Item(10, 14).Value = 485.325203
$ws.Cells.Item(10, 15).Value = 0.9790864123038654
$ws.Cells.Item(10, 16).Value = 0.9790864123038654
$ws.Cells.Item(10, 17).Value = 19.25101735224311
$ws.Cells.Item(10, 18).Value = 173.259156170188
$ws.Cells.Item(10, 19).Value = 0.002565987323455136
$ws.Cells.Item(10, 20).Value = 0.002565987323455136

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 0.1189986666666667
$ws.Cells.Item(11, 8).Value = 0.356996
$ws.Cells.Item(11, 9).Value = 0.002620797603979787
$ws.Cells.Item(11, 10).Value = 0.002620797603979787
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 0.67894
$ws.Cells.Item(11, 14).Value = 2.03682
$ws.Cells.Item(11, 15).Value = 0.004109044356199978
$ws.Cells.Item(11, 16).Value = 0.004109044356199979
$ws.Cells.Item(11, 17).Value = 0.08079295474666666
$ws.Cells.Item(11, 18).Value = 0.72713659272
$ws.Cells.Item(11, 19).Value = 0.00001076897360337557
$ws.Cells.Item(11, 20).Value = 0.00001076897360337557

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 0.1189986666666667
$ws.Cells.Item(12, 8).Value = 0.356996
$ws.Cells.Item(12, 9).Value = 0.002620797603979787
$ws.Cells.Item(12, 10).Value = 0.002620797603979787
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 1.763201333333333
$ws.Cells.Item(12, 14).Value = 5.289604
$ws.Cells.Item(12, 15).Value = 0.01067115280816804
$ws.Cells.Item(12, 16).Value = 0.01067115280816804
$ws.Cells.Item(12, 17).Value = 0.2098186077315555
$ws.Cells.Item(12, 18).Value = 1.888367469584
$ws.Cells.Item(12, 19).Value = 0.00002796693171134897
$ws.Cells.Item(12, 20).Value = 0.00002796693171134898

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 0.1189986666666667
$ws.Cells.Item(13, 8).Value = 0.356996
$ws.Cells.Item(13, 9).Value = 0.002620797603979787
$ws.Cells.Item(13, 10).Value = 0.002620797603979787
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 1.013424
$ws.Cells.Item(13, 14).Value = 3.040272
$ws.Cells.Item(13, 15).Value = 0.006133390531766587
$ws.Cells.Item(13, 16).Value = 0.006133390531766588
$ws.Cells.Item(13, 17).Value = 0.120596104768
$ws.Cells.Item(13, 18).Value = 1.085364942912
$ws.Cells.Item(13, 19).Value = 0.00001607437520992618
$ws.Cells.Item(13, 20).Value = 0.00001607437520992618

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 15.943029
$ws.Cells.Item(14, 8).Value = 47.829087
$ws.Cells.Item(14, 9).Value = 0.3511253812651704
$ws.Cells.Item(14, 10).Value = 0.3511253812651704
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 161.7750676666667
$ws.Cells.Item(14, 14).Value = 485.325203
$ws.Cells.Item(14, 15).Value = 0.9790864123038654
$ws.Cells.Item(14, 16).Value = 0.9790864123038654
$ws.Cells.Item(14, 17).Value = 2579.184595286629
$ws.Cells.Item(14, 18).Value = 23212.66135757966
$ws.Cells.Item(14, 19).Value = 0.3437820898117426
$ws.Cells.Item(14, 20).Value = 0.3437820898117426

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 15.943029
$ws.Cells.Item(15, 8).Value = 47.829087
$ws.Cells.Item(15, 9).Value = 0.3511253812651704
$ws.Cells.Item(15, 10).Value = 0.3511253812651704
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 0.67894
$ws.Cells.Item(15, 14).Value = 2.03682
$ws.Cells.Item(15, 15).Value = 0.004109044356199978
$ws.Cells.Item(15, 16).Value = 0.004109044356199979
$ws.Cells.Item(15, 17).Value = 10.82436010926
$ws.Cells.Item(15, 18).Value = 97.41924098334
$ws.Cells.Item(15, 19).Value = 0.001442789766206214
$ws.Cells.Item(15, 20).Value = 0.001442789766206214

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 15.943029
$ws.Cells.Item(16, 8).Value = 47.829087
$ws.Cells.Item(16, 9).Value = 0.3511253812651704
$ws.Cells.Item(16, 10).Value = 0.3511253812651704
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 1.763201333333333
$ws.Cells.Item(16, 14).Value = 5.289604
$ws.Cells.Item(16, 15).Value = 0.01067115280816804
$ws.Cells.Item(16, 16).Value = 0.01067115280816804
$ws.Cells.Item(16, 17).Value = 28.110769990172
$ws.Cells.Item(16, 18).Value = 252.996929911548
$ws.Cells.Item(16, 19).Value = 0.003746912598306897
$ws.Cells.Item(16, 20).Value = 0.003746912598306898

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 15.943029
$ws.Cells.Item(17, 8).Value = 47.829087
$ws.Cells.Item(17, 9).Value = 0.3511253812651704
$ws.Cells.Item(17, 10).Value = 0.3511253812651704
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 1.013424
$ws.Cells.Item(17, 14).Value = 3.040272
$ws.Cells.Item(17, 15).Value = 0.006133390531766587
$ws.Cells.Item(17, 16).Value = 0.006133390531766588
$ws.Cells.Item(17, 17).Value = 16.157048221296
$ws.Cells.Item(17, 18).Value = 145.413433991664
$ws.Cells.Item(17, 19).Value = 0.002153589088914729
$ws.Cells.Item(17, 20).Value = 0.002153589088914729
